$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 99; this pushes the existing rows 99-142 down to 100-143,
# preserving all of their data/formatting (date style included).
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new weekly record.
$ws.Range("A99").Value = 3
$ws.Range("B99").Value = "Femacal de La Calera"
$ws.Range("C99").Value = "Coquimbo"
$ws.Range("D99").Value = 44609
$ws.Range("E99").Value = 5
$ws.Range("F99").Value = 100112052
$ws.Range("G99").Value = "Albahaca"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 80
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = 4750
$ws.Range("N99").Value = "`$/docena de matas"
$ws.Range("O99").Value = "Provincia de Quillota"
$ws.Range("P99").Value = 792
$ws.Range("Q99").Value = 6
$ws.Range("R99").Value = "Hortaliza"

# Match the date number format used by the other rows in column D.
$ws.Range("D99").NumberFormat = $ws.Range("D100").NumberFormat
